# Updated cryptos list on Sat May 18 22:59:19 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.933.35"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "3.117.95"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "578.11"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").Value = "172.13"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("E9").Value = "  -3.30%  "

$ws.Range("E10").Value = "  -2.06%  "

$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").Value = "37.14"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").Value = "3.635.69"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "66.896.26"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "7.14"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "3.121.08"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "16.24"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").Value = "474.77"
$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("E22").Value = "  +4.68%  "

$ws.Range("D23").Value = "83.76"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "13.24"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("E25").Value = "  -3.80%  "

$ws.Range("D26").Value = "10.26"
$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("D28").Value = "7.88"
$ws.Range("E28").Value = "  -2.14%  "

$ws.Range("D29").Value = "2.37"
$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "28.59"
$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  -6.88%  "

# D34 ("0.999" -> "1.00") must keep the trailing zero as literal text, so
# force text interpretation with a leading apostrophe (same trick a human
# typing into Excel would use to stop "1.00" collapsing to the number 1).
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").Value = "0.975"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("D37").Value = "46.83"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").Value = "0.311"
$ws.Range("E40").Value = "  -2.66%  "

$ws.Range("E41").Value = "  +1.00%  "

# D42 ("8.61" -> "8.60") — same trailing-zero situation as D34 above.
$ws.Range("D42").Value = "'8.60"
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("D43").Value = "2.817.23"
$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").Value = "382.49"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").Value = "0.0352"
$ws.Range("E45").Value = "  -2.53%  "

$ws.Range("E46").Value = "  -9.92%  "

$ws.Range("D47").Value = "135.73"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "24.91"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("E51").Value = "  -0.90%  "
